# Commit: "updates with supplemental figures"
# Refresh the PSSM score grid (B2:K21) on Sheet1 with the recomputed
# supplemental-figure values. Header row (1) and the amino-acid labels in
# column A are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-16.97746999821544, 1.913426578242451, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544),
    @(-16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, 2.793479352539816, -16.97746999821544, -16.97746999821544),
    @(-16.97746999821544, 2.135827055261991, 2.186909114842111, -16.97746999821544, 3.409812168208096, -16.97746999821544, -16.97746999821544, -16.97746999821544, 2.264952851064117, -16.97746999821544),
    @(-16.97746999821544, 1.827521224529326, -16.97746999821544, -16.97746999821544, -16.97746999821544, 2.928188767594248, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544),
    @(-16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544),
    @(2.60629776041251, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544),
    @(-16.97746999821544, -16.97746999821544, -16.97746999821544, 1.841252592552709, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544),
    @(3.798099293384716, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544),
    @(-16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, 1.593204753560337, -16.97746999821544, 1.757494256590353),
    @(-16.97746999821544, -16.97746999821544, -16.97746999821544, 2.837846057495489, -16.97746999821544, 2.635846464752911, -16.97746999821544, -16.97746999821544, -16.97746999821544, 1.790316812006457),
    @(-16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544),
    @(-16.97746999821544, -16.97746999821544, -16.97746999821544, 2.440598776187593, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, 2.196160784519701, 2.04895888537495),
    @(-16.97746999821544, -16.97746999821544, 1.276400564095596, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, 2.196941277577776),
    @(-16.97746999821544, -16.97746999821544, 1.147585216104173, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544),
    @(-16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, 2.291686365335861, -16.97746999821544),
    @(-16.97746999821544, 2.07083859264506, 2.303058387631073, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, 1.16683473624505, 1.702288414521211, -16.97746999821544),
    @(-16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544, 1.532509849563959, 1.29893161202041, -16.97746999821544),
    @(-16.97746999821544, -16.97746999821544, 1.631019467900668, -16.97746999821544, -16.97746999821544, -16.97746999821544, 4.321917473725608, 1.441177368434952, -16.97746999821544, -16.97746999821544),
    @(-16.97746999821544, 0.9189392786280928, 1.47326201423671, -16.97746999821544, 3.228319672307146, -16.97746999821544, -16.97746999821544, 1.135015325135369, -16.97746999821544, 2.14949908827849),
    @(-16.97746999821544, 1.135320731752065, -16.97746999821544, 1.940534264069556, -16.97746999821544, 2.625955322627018, -16.97746999821544, -16.97746999821544, -16.97746999821544, -16.97746999821544)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 2).Value = $row[$j]
    }
}
